# Update the "umod title" column (A) values on the UMOD_TABLE sheet with the
# revised descriptive modification names.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = "acetylation of protein n-term"
    4  = "acetylation of K"
    6  = "amidation of peptide c-term"
    7  = "pyro-glu from n-term Q"
    10 = "carbamidomethyl C"
    11 = "carbamylation of n-term peptide"
    12 = "carbamylation of K"
    13 = "carboxymethyl C"
    17 = "deamidation of N"
    18 = "deamidation of N and Q"
    21 = "sulphone of M"
    25 = "formylation of protein n-term"
    26 = "formylation of peptide n-term"
    27 = "pyro-glu from n-term E"
    29 = "guanidination of K"
    30 = "ICAT light"
    31 = "ICAT heavy"
    41 = "iTRAQ117 on nterm"
    42 = "iTRAQ117 on K"
    43 = "iTRAQ117 on Y"
    44 = "iTRAQ8plex:13C(7)15N(1) on Y"
    45 = "iTRAQ8plex:13C(7)15N(1) on nterm"
    46 = "iTRAQ8plex:13C(7)15N(1) on K"
    47 = "O18 on peptide n-term"
    48 = "di-O18 on peptide n-term"
    49 = "homoserine"
    50 = "homoserine lactone"
    51 = "methylation of peptide c-term"
    52 = "methylation of D"
    53 = "methylation of E"
    54 = "MMTS on C"
    61 = "NIPCAM"
    62 = "oxidation of W"
    63 = "oxidation of H"
    64 = "oxidation of M"
    65 = "phosphorylation of Y"
    66 = "phosphorylation of T"
    67 = "phosphorylation of S"
    68 = "propionamide C"
    69 = "s-pyridylethylation of C"
    71 = "sulfation of Y"
    76 = "TMT duplex on n-term peptide"
    77 = "TMT duplex on K"
    78 = "TMT 6-plex on n-term peptide"
    79 = "TMT 6-plex on K"
}

foreach ($row in $updates.Keys) {
    $ws.Range("A$row").Value = $updates[$row]
}
